# Auto-generated script applying Phantom_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2998
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2998
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8994
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9330

# Row 19
$ws.Range("H19").Value = 1431.3334
$ws.Range("I19").Value = 1531.5714
$ws.Range("J19").Value = 1291
$ws.Range("K19").Value = 1531.5714
$ws.Range("L19").Value = 1291
$ws.Range("M19").Value = -1356.5714

# Row 33
$ws.Range("H33").Value = 3105.4443
$ws.Range("I33").Value = 2583.3333
$ws.Range("J33").Value = 4149.6665
$ws.Range("K33").Value = 2583.3333
$ws.Range("L33").Value = 4149.6665
$ws.Range("M33").Value = -2354.3333
$ws.Range("N33").Value = -4607.6665

# Row 58
$ws.Range("H58").Value = 4925.778
$ws.Range("I58").Value = 415
$ws.Range("J58").Value = 7181.1665
$ws.Range("K58").Value = 1245
$ws.Range("L58").Value = 21543.4995
$ws.Range("M58").Value = -1095
$ws.Range("N58").Value = -21843.4995

# Row 74
$ws.Range("H74").Value = 6552.1055
$ws.Range("I74").Value = 6552.1055
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6552.1055
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5616.1055
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 6552.1055
$ws.Range("I77").Value = 6552.1055
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 32760.5275
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -28080.5275
$ws.Range("N77").ClearContents()

# Row 107
$ws.Range("H107").Value = 1007.55554
$ws.Range("I107").Value = 796
$ws.Range("J107").Value = 2700
$ws.Range("K107").Value = 796
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 1124
$ws.Range("N107").Value = -6540

# Row 137
$ws.Range("H137").Value = 23811352
$ws.Range("I137").Value = 83334710
$ws.Range("J137").Value = 2008.9
$ws.Range("K137").Value = 250004130
$ws.Range("L137").Value = 6026.700000000001
$ws.Range("M137").Value = -250001580
$ws.Range("N137").Value = -11126.7

# Row 138
$ws.Range("H138").Value = 3475.0286
$ws.Range("I138").Value = 3904.8572
$ws.Range("J138").Value = 3188.476
$ws.Range("K138").Value = 11714.5716
$ws.Range("L138").Value = 9565.428
$ws.Range("M138").Value = -6574.571599999999
$ws.Range("N138").Value = -19845.428


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1459.7646
$ws.Range("I2").Value = 943.8333
$ws.Range("J2").Value = 2698
$ws.Range("K2").Value = 943.8333
$ws.Range("L2").Value = 2698
$ws.Range("M2").Value = -830.8333
$ws.Range("N2").Value = -2924

# Row 32
$ws.Range("H32").Value = 9144.909
$ws.Range("I32").Value = 5952.5884
$ws.Range("J32").Value = 19998.8
$ws.Range("K32").Value = 5952.5884
$ws.Range("L32").Value = 19998.8
$ws.Range("M32").Value = -5665.5884

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()

# Row 88
$ws.Range("H88").Value = 1187.5834
$ws.Range("I88").Value = 1162
$ws.Range("J88").Value = 1200.375
$ws.Range("K88").Value = 1162
$ws.Range("L88").Value = 1200.375
$ws.Range("M88").Value = -756
$ws.Range("N88").Value = -2012.375

# Row 91
$ws.Range("H91").Value = 1187.5834
$ws.Range("I91").Value = 1162
$ws.Range("J91").Value = 1200.375
$ws.Range("K91").Value = 1162
$ws.Range("L91").Value = 1200.375
$ws.Range("M91").Value = 242
$ws.Range("N91").Value = -4008.375

# Row 98
$ws.Range("H98").Value = 48999.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 48999.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 48999.5
$ws.Range("N98").Value = -54989.5

# Row 116
$ws.Range("H116").Value = 1459.7646
$ws.Range("I116").Value = 943.8333
$ws.Range("J116").Value = 2698
$ws.Range("K116").Value = 943.8333
$ws.Range("L116").Value = 2698
$ws.Range("M116").Value = 1350.1667
$ws.Range("N116").Value = -7286

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1459.7646
$ws.Range("I3").Value = 943.8333
$ws.Range("J3").Value = 2698
$ws.Range("K3").Value = 943.8333
$ws.Range("L3").Value = 2698
$ws.Range("M3").Value = -829.8333
$ws.Range("N3").Value = -2926


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 7328.125
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 7328.125
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 7328.125
$ws.Range("N28").Value = -7818.125

# Row 31
$ws.Range("H31").Value = 1453.6
$ws.Range("I31").Value = 1380.5555
$ws.Range("J31").Value = 2111
$ws.Range("K31").Value = 1380.5555
$ws.Range("L31").Value = 2111
$ws.Range("M31").Value = -1085.5555

# Row 34
$ws.Range("H34").Value = 1453.6
$ws.Range("I34").Value = 1380.5555
$ws.Range("J34").Value = 2111
$ws.Range("K34").Value = 1380.5555
$ws.Range("L34").Value = 2111
$ws.Range("M34").Value = -1178.5555

# Row 93
$ws.Range("H93").Value = 21802.334
$ws.Range("I93").Value = 21802.334
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 21802.334
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -19930.334

# Row 141
$ws.Range("H141").Value = 375614.62
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 422131
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 422131
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -432491


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 111.53333
$ws.Range("I23").Value = 32.666668
$ws.Range("J23").Value = 164.11111
$ws.Range("K23").Value = 98.000004
$ws.Range("L23").Value = 492.33333
$ws.Range("M23").Value = 136.999996
$ws.Range("N23").Value = -962.3333299999999

# Row 24
$ws.Range("H24").Value = 2874.75
$ws.Range("I24").Value = 2500
$ws.Range("J24").Value = 3249.5
$ws.Range("K24").Value = 7500
$ws.Range("L24").Value = 9748.5
$ws.Range("M24").Value = -7270
$ws.Range("N24").Value = -10208.5

# Row 34
$ws.Range("H34").Value = 37989.242
$ws.Range("I34").Value = 524.5
$ws.Range("J34").Value = 40764.406
$ws.Range("K34").Value = 1573.5
$ws.Range("L34").Value = 122293.218
$ws.Range("M34").Value = -1489.5
$ws.Range("N34").Value = -122461.218

# Row 38
$ws.Range("H38").Value = 53.692307
$ws.Range("I38").Value = 60.8
$ws.Range("J38").Value = 49.25
$ws.Range("K38").Value = 182.4
$ws.Range("L38").Value = 147.75
$ws.Range("M38").Value = 164.6
$ws.Range("N38").Value = -841.75

# Row 58
$ws.Range("H58").Value = 5749.5
$ws.Range("I58").Value = 5749.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 17248.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -17120.5
$ws.Range("N58").ClearContents()

# Row 129
$ws.Range("H129").Value = 1606
$ws.Range("I129").Value = 993.6667
$ws.Range("J129").Value = 2524.5
$ws.Range("K129").Value = 2981.0001
$ws.Range("L129").Value = 7573.5
$ws.Range("M129").Value = 2018.9999
$ws.Range("N129").Value = -17573.5

# Row 131
$ws.Range("H131").Value = 2347.0667
$ws.Range("I131").Value = 2206.5557
$ws.Range("J131").Value = 2557.8333
$ws.Range("K131").Value = 6619.6671
$ws.Range("L131").Value = 7673.499899999999
$ws.Range("M131").Value = -1579.6671
$ws.Range("N131").Value = -17753.4999


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2042.9286
$ws.Range("I102").Value = 2042.9286
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2042.9286
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -420.9286

# Row 105
$ws.Range("H105").Value = 29332
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 29332
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 29332
$ws.Range("N105").Value = -36320

# Row 132
$ws.Range("I132").Value = 4261.25
$ws.Range("J132").Value = 55559424
$ws.Range("K132").Value = 12783.75
$ws.Range("L132").Value = 166678272
$ws.Range("M132").Value = -10253.75


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 50001870
$ws.Range("I22").Value = 1057
$ws.Range("J22").Value = 83335750
$ws.Range("K22").Value = 1057
$ws.Range("L22").Value = 83335750
$ws.Range("M22").Value = -762
$ws.Range("N22").Value = -83336340

# Row 27
$ws.Range("H27").Value = 50001870
$ws.Range("I27").Value = 1057
$ws.Range("J27").Value = 83335750
$ws.Range("K27").Value = 1057
$ws.Range("L27").Value = 83335750
$ws.Range("M27").Value = -950
$ws.Range("N27").Value = -83335964

# Row 94
$ws.Range("H94").Value = 82999.664
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 82999.664
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 82999.664
$ws.Range("N94").Value = -84351.664

# Row 106
$ws.Range("H106").Value = 18576.666
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 18576.666
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 18576.666
$ws.Range("N106").Value = -21100.666

# Row 112
$ws.Range("H112").Value = 32000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 32000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 32000
$ws.Range("N112").Value = -34954

# Row 132
$ws.Range("H132").Value = 3680.9375
$ws.Range("I132").Value = 3949.5
$ws.Range("J132").Value = 3591.4167
$ws.Range("K132").Value = 11848.5
$ws.Range("L132").Value = 10774.2501
$ws.Range("M132").Value = -9318.5
$ws.Range("N132").Value = -15834.2501


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 80000
$ws.Range("I8").Value = 80000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 80000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -79860
$ws.Range("N8").ClearContents()

# Row 15
$ws.Range("H15").Value = 7500
$ws.Range("I15").Value = 7500
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7500
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -7212

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 97
$ws.Range("H97").Value = 19999
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 19999
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 19999
$ws.Range("N97").Value = -21981

# Row 107
$ws.Range("H107").Value = 4996.5
$ws.Range("I107").Value = 4996.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 14989.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -13069.5
$ws.Range("N107").ClearContents()

# Row 124
$ws.Range("H124").Value = 25685.2
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 25685.2
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 25685.2
$ws.Range("N124").Value = -35505.2

# Row 136
$ws.Range("H136").Value = 7050.579
$ws.Range("I136").Value = 8195.625
$ws.Range("J136").Value = 943.6667
$ws.Range("K136").Value = 24586.875
$ws.Range("L136").Value = 2831.0001
$ws.Range("M136").Value = -22036.875
$ws.Range("N136").Value = -7931.0001

